# SSHConfig.xlsx: point the SSH target at a different experiment machine.
#  - ipaddress (B2): 192.168.122.1  -> 130.237.10.123
#  - hostname  (B7): adeye06u       -> adeye03u
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "130.237.10.123"
$ws.Range("B7").Value = "adeye03u"

# Leave the cursor on the cell that was last edited, matching the saved
# selection state in the workbook.
$ws.Range("B7").Select()
